$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value2 = 8.75
$ws.Range("I11").Value2 = 8.75
$ws.Range("K11").Value2 = 8.75
$ws.Range("M11").Value2 = 131.25
$ws.Range("H15").Value2 = 122.5
$ws.Range("I15").Value2 = 122.5
$ws.Range("K15").Value2 = 367.5
$ws.Range("M15").Value2 = -198.5
$ws.Range("H33").Value2 = 963.8570999999999
$ws.Range("I33").Value2 = 349.6
$ws.Range("J33").Value2 = 2499.5
$ws.Range("K33").Value2 = 349.6
$ws.Range("L33").Value2 = 2499.5
$ws.Range("M33").Value2 = -120.6
$ws.Range("N33").Value2 = -2957.5
$ws.Range("H39").Value2 = 574.1111
$ws.Range("I39").Value2 = 520.875
$ws.Range("J39").Value2 = 1000
$ws.Range("K39").Value2 = 1562.625
$ws.Range("L39").Value2 = 3000
$ws.Range("M39").Value2 = -1266.625
$ws.Range("N39").Value2 = -3592
$ws.Range("H58").Value2 = 1362.5
$ws.Range("I58").Value2 = 1125
$ws.Range("J58").Value2 = 2075
$ws.Range("K58").Value2 = 3375
$ws.Range("L58").Value2 = 6225
$ws.Range("M58").Value2 = -3225
$ws.Range("N58").Value2 = -6525
$ws.Range("H98").Value2 = 698.2857
$ws.Range("I98").Value2 = 698.2857
$ws.Range("K98").Value2 = 698.2857
$ws.Range("M98").Value2 = 799.7143
$ws.Range("H106").Value2 = 4999
$ws.Range("I106").Value2 = 4999
$ws.Range("K106").Value2 = 4999
$ws.Range("M106").Value2 = -4368
$ws.Range("H122").Value2 = 698.2857
$ws.Range("I122").Value2 = 698.2857
$ws.Range("K122").Value2 = 2094.8571
$ws.Range("M122").Value2 = 355.1428999999998
$ws.Range("H129").Value2 = 533.3333
$ws.Range("J129").Value2 = 0
$ws.Range("L129").Value2 = 0
$ws.Range("N129").ClearContents()
$ws.Range("H138").Value2 = 2763.0908
$ws.Range("J138").Value2 = 3060
$ws.Range("L138").Value2 = 9180
$ws.Range("N138").Value2 = -19460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1642.0714
$ws.Range("I61").Value2 = 1642.0714
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 1642.0714
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -1430.0714
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value2 = 1000
$ws.Range("I74").Value2 = 1000
$ws.Range("K74").Value2 = 1000
$ws.Range("M74").Value2 = -126
$ws.Range("H77").Value2 = 1000
$ws.Range("I77").Value2 = 1000
$ws.Range("K77").Value2 = 5000
$ws.Range("M77").Value2 = -632
$ws.Range("H88").Value2 = 2262.75
$ws.Range("I88").Value2 = 1809.6
$ws.Range("K88").Value2 = 1809.6
$ws.Range("M88").Value2 = -1403.6
$ws.Range("H91").Value2 = 2262.75
$ws.Range("I91").Value2 = 1809.6
$ws.Range("K91").Value2 = 1809.6
$ws.Range("M91").Value2 = -405.5999999999999
$ws.Range("H97").Value2 = 2213.6
$ws.Range("I97").Value2 = 1861.8462
$ws.Range("K97").Value2 = 1861.8462
$ws.Range("M97").Value2 = -1365.8462
$ws.Range("H102").Value2 = 3250
$ws.Range("I102").Value2 = 3250
$ws.Range("K102").Value2 = 3250
$ws.Range("M102").Value2 = -1628
$ws.Range("H136").Value2 = 1642.0714
$ws.Range("I136").Value2 = 1642.0714
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 4926.2142
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -2376.2142
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value2 = 80
$ws.Range("I8").Value2 = 80
$ws.Range("J8").Value2 = 0
$ws.Range("K8").Value2 = 80
$ws.Range("L8").Value2 = 0
$ws.Range("M8").Value2 = 60
$ws.Range("N8").ClearContents()
$ws.Range("H20").Value2 = 1601.8
$ws.Range("I20").Value2 = 1000
$ws.Range("K20").Value2 = 1000
$ws.Range("M20").Value2 = -753
$ws.Range("H86").Value2 = 2782.8333
$ws.Range("I86").Value2 = 4266.3335
$ws.Range("J86").Value2 = 1299.3334
$ws.Range("K86").Value2 = 4266.3335
$ws.Range("L86").Value2 = 1299.3334
$ws.Range("M86").Value2 = -3143.3335
$ws.Range("N86").Value2 = -3545.3334
$ws.Range("H89").Value2 = 2782.8333
$ws.Range("I89").Value2 = 4266.3335
$ws.Range("J89").Value2 = 1299.3334
$ws.Range("K89").Value2 = 21331.6675
$ws.Range("L89").Value2 = 6496.666999999999
$ws.Range("M89").Value2 = -15715.6675
$ws.Range("N89").Value2 = -17728.667
$ws.Range("H94").Value2 = 2845.3333
$ws.Range("I94").Value2 = 2683.2222
$ws.Range("J94").Value2 = 3331.6667
$ws.Range("K94").Value2 = 2683.2222
$ws.Range("L94").Value2 = 3331.6667
$ws.Range("M94").Value2 = -2232.2222
$ws.Range("N94").Value2 = -4233.6667
$ws.Range("H105").Value2 = 3971.4
$ws.Range("I105").Value2 = 3089.25
$ws.Range("K105").Value2 = 3089.25
$ws.Range("M105").Value2 = -1342.25
$ws.Range("H107").Value2 = 758.5
$ws.Range("I107").Value2 = 595.4286
$ws.Range("K107").Value2 = 595.4286
$ws.Range("M107").Value2 = 1324.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value2 = 797.5
$ws.Range("I19").Value2 = 797.5
$ws.Range("K19").Value2 = 797.5
$ws.Range("M19").Value2 = -627.5
$ws.Range("H24").Value2 = 797.5
$ws.Range("I24").Value2 = 797.5
$ws.Range("K24").Value2 = 797.5
$ws.Range("M24").Value2 = -627.5
$ws.Range("H31").Value2 = 2255.3333
$ws.Range("I31").Value2 = 1599.6
$ws.Range("K31").Value2 = 1599.6
$ws.Range("M31").Value2 = -1304.6
$ws.Range("H34").Value2 = 2255.3333
$ws.Range("I34").Value2 = 1599.6
$ws.Range("K34").Value2 = 1599.6
$ws.Range("M34").Value2 = -1397.6
$ws.Range("H58").Value2 = 736.5714
$ws.Range("I58").Value2 = 519.4545000000001
$ws.Range("K58").Value2 = 519.4545000000001
$ws.Range("M58").Value2 = -316.4545000000001
$ws.Range("H99").Value2 = 5045.5
$ws.Range("I99").Value2 = 5045.5
$ws.Range("K99").Value2 = 5045.5
$ws.Range("M99").Value2 = -3547.5
$ws.Range("H122").Value2 = 1290.2858
$ws.Range("I122").Value2 = 1338.6666
$ws.Range("K122").Value2 = 4015.9998
$ws.Range("M122").Value2 = -1565.9998
$ws.Range("H126").Value2 = 5045.5
$ws.Range("I126").Value2 = 5045.5
$ws.Range("K126").Value2 = 15136.5
$ws.Range("M126").Value2 = -12666.5
$ws.Range("H132").Value2 = 1798.7587
$ws.Range("I132").Value2 = 1698.6154
$ws.Range("K132").Value2 = 5095.8462
$ws.Range("M132").Value2 = -2565.8462
$ws.Range("H134").Value2 = 2030.6666
$ws.Range("I134").Value2 = 1941.2727
$ws.Range("K134").Value2 = 5823.8181
$ws.Range("M134").Value2 = -3288.8181
$ws.Range("H136").Value2 = 736.5714
$ws.Range("I136").Value2 = 519.4545000000001
$ws.Range("K136").Value2 = 1558.3635
$ws.Range("M136").Value2 = 991.6364999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value2 = 349.33334
$ws.Range("I44").Value2 = 474
$ws.Range("J44").Value2 = 100
$ws.Range("K44").Value2 = 1422
$ws.Range("L44").Value2 = 300
$ws.Range("M44").Value2 = -1024
$ws.Range("N44").Value2 = -1096
$ws.Range("H103").Value2 = 36188.145
$ws.Range("J103").Value2 = 950
$ws.Range("L103").Value2 = 2850
$ws.Range("N103").Value2 = -4608
$ws.Range("H131").Value2 = 1317.4
$ws.Range("I131").Value2 = 499.66666
$ws.Range("J131").Value2 = 1461.7059
$ws.Range("K131").Value2 = 1498.99998
$ws.Range("L131").Value2 = 4385.1177
$ws.Range("M131").Value2 = 3541.00002
$ws.Range("N131").Value2 = -14465.1177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value2 = 500
$ws.Range("I29").Value2 = 500
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 500
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = -210
$ws.Range("N29").ClearContents()
$ws.Range("H44").Value2 = 17500.5
$ws.Range("J44").Value2 = 17500.5
$ws.Range("L44").Value2 = 17500.5
$ws.Range("N44").Value2 = -18692.5
$ws.Range("H46").Value2 = 8456.9
$ws.Range("I46").Value2 = 1526.3334
$ws.Range("J46").Value2 = 11427.143
$ws.Range("K46").Value2 = 1526.3334
$ws.Range("L46").Value2 = 11427.143
$ws.Range("M46").Value2 = -1370.3334
$ws.Range("N46").Value2 = -11739.143
$ws.Range("H55").Value2 = 1999
$ws.Range("I55").Value2 = 1999
$ws.Range("K55").Value2 = 1999
$ws.Range("M55").Value2 = -1672
$ws.Range("H97").Value2 = 2382.25
$ws.Range("I97").Value2 = 2293.1667
$ws.Range("K97").Value2 = 2293.1667
$ws.Range("M97").Value2 = -1797.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 8485.571
$ws.Range("J22").Value2 = 9250
$ws.Range("L22").Value2 = 9250
$ws.Range("N22").Value2 = -9840
$ws.Range("H27").Value2 = 8485.571
$ws.Range("J27").Value2 = 9250
$ws.Range("L27").Value2 = 9250
$ws.Range("N27").Value2 = -9464
$ws.Range("H32").Value2 = 2316.5
$ws.Range("I32").Value2 = 2316.5
$ws.Range("K32").Value2 = 2316.5
$ws.Range("M32").Value2 = -1999.5
$ws.Range("H93").Value2 = 591.1429000000001
$ws.Range("I93").Value2 = 591.1429000000001
$ws.Range("K93").Value2 = 591.1429000000001
$ws.Range("M93").Value2 = 656.8570999999999
$ws.Range("H110").Value2 = 0
$ws.Range("J110").Value2 = 0
$ws.Range("L110").Value2 = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("N19").ClearContents()
$ws.Range("H113").Value2 = 1338
$ws.Range("I113").Value2 = 1341.1666
$ws.Range("K113").Value2 = 4023.4998
$ws.Range("M113").Value2 = -1853.4998
